# Lunggo_Config.xlsx - "Adding Mandiri Account in Lunggo Config"
#
# Adds four new config rows (rows 59-62) for the "mandiri" section:
#   mandiri.webCompanyId      = TMDZ001
#   mandiri.webUserName       = rama_maker_1
#   mandiri.webPassword       = 164926a78b265daf9cb7c15dcbbe2a5b6f5074ae
#   mandiri.bankAccountNumber = 1020006675802 (numeric, left-aligned integer format)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: "*" for rows 59-62 ---
$ws.Range("A59:A62").Value = "*"

# --- Column B: "mandiri" for rows 59-62 ---
$ws.Range("B59:B62").Value = "mandiri"

# --- Column C: per-row keys (webUserName / webPassword already exist as shared strings) ---
$ws.Range("C59").Value = "webCompanyId"
$ws.Range("C60").Value = "webUserName"
$ws.Range("C61").Value = "webPassword"
$ws.Range("C62").Value = "bankAccountNumber"

# --- Column D: formulas producing "@@.*.mandiri.<key>@@" ---
$ws.Range("D59").Formula = '="@@."&A59&"."&B59&"."&C59&"@@"'
$ws.Range("D60").Formula = '="@@."&A60&"."&B60&"."&C60&"@@"'
$ws.Range("D61").Formula = '="@@."&A61&"."&B61&"."&C61&"@@"'
$ws.Range("D62").Formula = '="@@."&A62&"."&B62&"."&C62&"@@"'

# --- Columns E:H values for rows 59-61 (same value repeated across all 4 environments) ---
$ws.Range("E59:H59").Value = "TMDZ001"
$ws.Range("E59:H59").WrapText = $true

$ws.Range("E60:H60").Value = "rama_maker_1"
$ws.Range("E60:H60").WrapText = $true

$ws.Range("E61:H61").Value = "164926a78b265daf9cb7c15dcbbe2a5b6f5074ae"
$ws.Range("E61:H61").WrapText = $true

# --- Row 62: numeric bank account number, with a left-aligned integer number format ---
# Build the exact target format on a scratch cell first (inheriting border/wrap/vertical
# alignment from an already-formatted cell) so the workbook doesn't fragment its style
# table with one-off intermediate styles, then stamp that format onto E62:H62.
$ws.Range("E58").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$helper = $ws.Range("Z1")
$helper.HorizontalAlignment = -4131
$helper.NumberFormat = "0"
$helper.Copy()

$dest = $ws.Range("E62:H62")
$dest.PasteSpecial(-4122)
$dest.Value = 1020006675802
$ws.Range("Z1").Clear()

# --- Update the visible selection to match the edited area ---
$ws.Range("C60").Select()
